$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44292
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 14000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 967

# Row 3
$ws.Range("D3").Value = 44301
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 17500
$ws.Range("O3").Value = 17500
$ws.Range("P3").Value = 17500
$ws.Range("S3").Value = 1167

# Row 4
$ws.Range("D4").Value = 44301
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 14000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 967

# Row 5
$ws.Range("D5").Value = 44298
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 967

# Row 6
$ws.Range("D6").Value = 44348
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 1000

# Row 7
$ws.Range("D7").Value = 44348
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13500
$ws.Range("P7").Value = 13250
$ws.Range("S7").Value = 883

# Row 8
$ws.Range("D8").Value = 44305
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 17500
$ws.Range("O8").Value = 17500
$ws.Range("P8").Value = 17500
$ws.Range("S8").Value = 1167

# Row 9
$ws.Range("D9").Value = 44305
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("S9").Value = 967

# Row 10
$ws.Range("D10").Value = 44351
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 100
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("S10").Value = 1000

# Row 11
$ws.Range("D11").Value = 44351
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 13000
$ws.Range("O11").Value = 13500
$ws.Range("P11").Value = 13250
$ws.Range("S11").Value = 883

# Row 12
$ws.Range("D12").Value = 44293
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 60
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 14500
$ws.Range("S12").Value = 967

# Row 13
$ws.Range("D13").Value = 44336
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 17000
$ws.Range("O13").Value = 17000
$ws.Range("P13").Value = 17000
$ws.Range("S13").Value = 1133

# Row 14
$ws.Range("D14").Value = 44336
$ws.Range("M14").Value = 120
$ws.Range("O14").Value = 14500
$ws.Range("P14").Value = 14250
$ws.Range("S14").Value = 950

# Row 15
$ws.Range("D15").Value = 44327
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 17000
$ws.Range("P15").Value = 17000
$ws.Range("S15").Value = 1133

# Row 16
$ws.Range("D16").Value = 44327

# Row 17
$ws.Range("D17").Value = 44309
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 17500
$ws.Range("O17").Value = 17500
$ws.Range("P17").Value = 17500
$ws.Range("S17").Value = 1167

# Row 18
$ws.Range("D18").Value = 44309
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 14500
$ws.Range("P18").Value = 14250
$ws.Range("S18").Value = 950

# Row 19
$ws.Range("D19").Value = 44295
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 160

# Row 20
$ws.Range("D20").Value = 44313
$ws.Range("L20").Value = "Especial"
$ws.Range("Q20").Value = "$/caja 14 kilos empedrada"
$ws.Range("S20").Value = 1250
$ws.Range("T20").Value = 14

# Row 21
$ws.Range("D21").Value = 44313
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 16000
$ws.Range("O21").Value = 16000
$ws.Range("P21").Value = 16000
$ws.Range("Q21").Value = "$/caja 14 kilos empedrada"
$ws.Range("S21").Value = 1143
$ws.Range("T21").Value = 14

# Row 22
$ws.Range("D22").Value = 44313
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 80
$ws.Range("O22").Value = 14000
$ws.Range("P22").Value = 14000
$ws.Range("Q22").Value = "$/caja 14 kilos empedrada"
$ws.Range("S22").Value = 1000
$ws.Range("T22").Value = 14

# Row 23
$ws.Range("D23").Value = 44344
$ws.Range("N23").Value = 16000
$ws.Range("O23").Value = 16000
$ws.Range("P23").Value = 16000
$ws.Range("S23").Value = 1067

# Row 24
$ws.Range("D24").Value = 44344
$ws.Range("M24").Value = 120
$ws.Range("N24").Value = 13000
$ws.Range("O24").Value = 13500
$ws.Range("P24").Value = 13250
$ws.Range("S24").Value = 883

# Row 25
$ws.Range("D25").Value = 44334
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 15500
$ws.Range("Q25").Value = "$/caja 15 kilos empedrada"
$ws.Range("S25").Value = 1033
$ws.Range("T25").Value = 15

# Row 26
$ws.Range("D26").Value = 44334
$ws.Range("L26").Value = "Segunda"
$ws.Range("N26").Value = 14500
$ws.Range("O26").Value = 14500
$ws.Range("P26").Value = 14500
$ws.Range("Q26").Value = "$/caja 15 kilos empedrada"
$ws.Range("S26").Value = 967
$ws.Range("T26").Value = 15

# Row 27
$ws.Range("D27").Value = 44330
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 17000
$ws.Range("O27").Value = 17000
$ws.Range("P27").Value = 17000
$ws.Range("Q27").Value = "$/caja 15 kilos empedrada"
$ws.Range("S27").Value = 1133
$ws.Range("T27").Value = 15

# Row 28
$ws.Range("D28").Value = 44330
$ws.Range("L28").Value = "Segunda"
$ws.Range("O28").Value = 14500
$ws.Range("P28").Value = 14250
$ws.Range("S28").Value = 950

# Row 29
$ws.Range("D29").Value = 44285
$ws.Range("M29").Value = 160
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 15500
$ws.Range("S29").Value = 1033

# Row 30
$ws.Range("D30").Value = 44316
$ws.Range("N30").Value = 17500
$ws.Range("O30").Value = 17500
$ws.Range("P30").Value = 17500
$ws.Range("S30").Value = 1167

# Row 31
$ws.Range("D31").Value = 44316

# Row 32
$ws.Range("D32").Value = 44323
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 17000
$ws.Range("O32").Value = 17000
$ws.Range("P32").Value = 17000
$ws.Range("S32").Value = 1133

# Row 33
$ws.Range("D33").Value = 44323
$ws.Range("L33").Value = "Segunda"
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 14000
$ws.Range("O33").Value = 14000
$ws.Range("P33").Value = 14000
$ws.Range("S33").Value = 933

# Row 34
$ws.Range("D34").Value = 44306
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 100
$ws.Range("N34").Value = 17500
$ws.Range("O34").Value = 17500
$ws.Range("P34").Value = 17500
$ws.Range("S34").Value = 1167

# Row 35
$ws.Range("D35").Value = 44306
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 200
$ws.Range("N35").Value = 14000
$ws.Range("O35").Value = 14500
$ws.Range("P35").Value = 14250
$ws.Range("S35").Value = 950

# Row 36
$ws.Range("D36").Value = 44299
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 60
$ws.Range("N36").Value = 17500
$ws.Range("O36").Value = 17500
$ws.Range("P36").Value = 17500
$ws.Range("S36").Value = 1167

# Row 37
$ws.Range("D37").Value = 44299
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 120

# Row 38
$ws.Range("D38").Value = 44302
$ws.Range("N38").Value = 17500
$ws.Range("O38").Value = 17500
$ws.Range("P38").Value = 17500
$ws.Range("S38").Value = 1167

# Row 39
$ws.Range("D39").Value = 44302
$ws.Range("O39").Value = 15000
$ws.Range("P39").Value = 14500
$ws.Range("S39").Value = 967
